# Applies the "Built site for gh-pages" styles.xml template refresh:
#   1. Adds a new "AbstractTitle" paragraph style (based on Normal,
#      followed by Abstract).
#   2. Tightens the "Abstract" style's space-before from 15pt to 5pt.
#   3. Adds a new built-in-style "FootnoteBlockText" paragraph style
#      (based on Footnote Text, followed by Footnote Text).

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" style ---------------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060   # OLE BGR for RGB 345A8A

# --- 2. Tighten "Abstract" style's space-before (300 -> 100 twips) -------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. New "Footnote Block Text" style -----------------------------------
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "Styles updated: AbstractTitle added, Abstract before-spacing tightened, FootnoteBlockText added."
